$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Numeric-looking text values (Price column D, Hora column G) ---
# Force each cell to Text format individually (union ranges only apply to the first area)
# so these remain strings, matching the original inlineStr text values.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "245.59"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "17"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "21.96"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "17"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.367"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "17"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05855"
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "17"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.395"
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "17"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.360"
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "17"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8156"
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "17"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.016"
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "17"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.01126"
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "17"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1421"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "17"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.04056"
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "17"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07394"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "17"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.02971"
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "17"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.145"
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "17"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.09393"
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "17"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.001594"
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = "17"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.04817"
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = "17"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006072"
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "17"
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = "17"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0009844"
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = "17"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0001500"
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = "17"
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = "17"
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = "17"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.3239"
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = "17"
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = "17"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0002484"
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = "17"
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = "17"
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = "17"
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = "17"
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = "17"
$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = "17"
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = "17"
$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = "17"
$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = "17"
$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = "17"
$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = "17"
$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = "17"
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = "17"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03859"
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = "17"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006372"
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = "17"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1073"
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = "17"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002601"
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = "17"
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = "17"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005631"
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = "17"
$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = "17"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.7702"
$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = "17"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.08753"
$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = "17"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002101"
$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = "17"
$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = "17"
$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = "17"

# --- Text values (Coin name B, Link C, Volume(1h) E) ---
$ws.Range("B10").Value = "One"
$ws.Range("C10").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("E10").Value = "9OneONEBestin24h"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("E11").Value = "10WazirXWRX"
$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("B13").Value = "MandalaExchangeToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("E13").Value = "12MandalaExchangeTokenMDX"
$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("E14").Value = "13BitrueCoinBTR"
$ws.Range("B15").Value = "MCDex"
$ws.Range("C15").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("E15").Value = "14MCDexMCB"
$ws.Range("B16").Value = "BitMartToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("E16").Value = "15BitMartTokenBMX"
$ws.Range("B17").Value = "BitForexToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("E17").Value = "16BitForexTokenBF"
$ws.Range("B18").Value = "CoinExToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("E18").Value = "17CoinExTokenCET"
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("E42").Value = "41BKEXTokenBKK"
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("E43").Value = "42CEJICEJI"
$ws.Range("E48").Value = "47BOLOBOLOWorstin24h"

Write-Host "Applied all changes"
